$wb = $excel.ActiveWorkbook

# --- PIR sheet: rows 224-239, motion sensor "No Motion"/"Inactive" ---
$wsPIR = $wb.Worksheets.Item("PIR")
$pirData = @(
    ,@(224, "12:21:33")
    ,@(225, "12:21:34")
    ,@(226, "12:21:37")
    ,@(227, "12:21:38")
    ,@(228, "12:21:39")
    ,@(229, "12:21:41")
    ,@(230, "12:21:44")
    ,@(231, "12:21:45")
    ,@(232, "12:21:54")
    ,@(233, "12:21:58")
    ,@(234, "12:22:03")
    ,@(235, "12:22:08")
    ,@(236, "12:22:14")
    ,@(237, "12:22:18")
    ,@(238, "12:22:23")
    ,@(239, "12:22:28")
)
$pirFirstRow = 224
$pirLastRow = 239
$rngPirA = $wsPIR.Range("A" + $pirFirstRow + ":A" + $pirLastRow)
$rngPirA.NumberFormat = "@"
foreach ($row in $pirData) {
    $r = $row[0]
    $t = $row[1]
    $wsPIR.Cells.Item($r, 1).Value = "2026-01-28"
    $wsPIR.Cells.Item($r, 2).Value = $t
    $wsPIR.Cells.Item($r, 3).Value = "12:00"
    $wsPIR.Cells.Item($r, 4).Value = "Bathroom"
    $wsPIR.Cells.Item($r, 5).Value = "No Motion"
    $wsPIR.Cells.Item($r, 6).Value = "Inactive"
}
$rngPirA.ClearFormats()

# --- Humidity sheet: rows 211-225, "Active" with % values ---
$wsHum = $wb.Worksheets.Item("Humidity")
$humidityData = @(
    ,@(211, "12:21:32", "86.7%")
    ,@(212, "12:21:34", "87.6%")
    ,@(213, "12:21:36", "86.7%")
    ,@(214, "12:21:38", "86.7%")
    ,@(215, "12:21:40", "87.7%")
    ,@(216, "12:21:42", "87.7%")
    ,@(217, "12:21:43", "86.7%")
    ,@(218, "12:21:46", "86.7%")
    ,@(219, "12:21:49", "87.4%")
    ,@(220, "12:21:53", "86.5%")
    ,@(221, "12:22:05", "85.9%")
    ,@(222, "12:22:13", "86.4%")
    ,@(223, "12:22:17", "87.4%")
    ,@(224, "12:22:25", "86.4%")
    ,@(225, "12:22:29", "87.3%")
)
$humFirstRow = 211
$humLastRow = 225
$rngHumA = $wsHum.Range("A" + $humFirstRow + ":A" + $humLastRow)
$rngHumA.NumberFormat = "@"
$rngHumE = $wsHum.Range("E" + $humFirstRow + ":E" + $humLastRow)
$rngHumE.NumberFormat = "@"
foreach ($row in $humidityData) {
    $r = $row[0]
    $t = $row[1]
    $v = $row[2]
    $wsHum.Cells.Item($r, 1).Value = "2026-01-28"
    $wsHum.Cells.Item($r, 2).Value = $t
    $wsHum.Cells.Item($r, 3).Value = "12:00"
    $wsHum.Cells.Item($r, 4).Value = "Bathroom"
    $wsHum.Cells.Item($r, 5).Value = $v
    $wsHum.Cells.Item($r, 6).Value = "Active"
}
$rngHumA.ClearFormats()
$rngHumE.ClearFormats()

# --- Temperature sheet: rows 211-225, "Active" with Celsius values ---
$wsTemp = $wb.Worksheets.Item("Temperature")
$temperatureData = @(
    ,@(211, "12:21:33", "23.0C")
    ,@(212, "12:21:35", "23.0C")
    ,@(213, "12:21:36", "23.0C")
    ,@(214, "12:21:39", "23.0C")
    ,@(215, "12:21:41", "23.0C")
    ,@(216, "12:21:43", "23.0C")
    ,@(217, "12:21:44", "23.0C")
    ,@(218, "12:21:46", "23.0C")
    ,@(219, "12:21:49", "23.1C")
    ,@(220, "12:21:54", "23.1C")
    ,@(221, "12:22:06", "23.1C")
    ,@(222, "12:22:14", "23.1C")
    ,@(223, "12:22:18", "23.1C")
    ,@(224, "12:22:26", "23.0C")
    ,@(225, "12:22:30", "23.1C")
)
$tempFirstRow = 211
$tempLastRow = 225
$rngTempA = $wsTemp.Range("A" + $tempFirstRow + ":A" + $tempLastRow)
$rngTempA.NumberFormat = "@"
foreach ($row in $temperatureData) {
    $r = $row[0]
    $t = $row[1]
    $v = $row[2]
    $wsTemp.Cells.Item($r, 1).Value = "2026-01-28"
    $wsTemp.Cells.Item($r, 2).Value = $t
    $wsTemp.Cells.Item($r, 3).Value = "12:00"
    $wsTemp.Cells.Item($r, 4).Value = "Bathroom"
    $wsTemp.Cells.Item($r, 5).Value = $v
    $wsTemp.Cells.Item($r, 6).Value = "Active"
}
$rngTempA.ClearFormats()

